$wb = $excel.ActiveWorkbook

# --- 1. Text change: "Ready for handoff" -> "In Translation" ---
# This status string appears in the Overview sheet (columns E & F, the
# per-language status columns) and in each language sheet's "Status" column.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- 2. Column width changes (report columns narrowed) ---
# Target OOXML stored width is ~13.41 characters; the host's ColumnWidth
# setter quantizes to whole-pixel steps, so 12.5 lands on the closest
# reachable stored width (13.33).
# Overview: columns E and F (zh-cn / de-de status columns)
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

# zh-cn / de-de: column C (Status column)
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
